# analyze 120524 and improve algorithm
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 29
$ws.Range("C4").Value = 9
$ws.Range("C5").Value = 9
$ws.Range("C12").Value = 65
$ws.Range("C18").Value = 20
$ws.Range("C21").Value = 27
$ws.Range("C22").Value = 4
$ws.Range("C24").Value = 0
$ws.Range("C29").Value = 94
$ws.Range("C30").Value = 12
$ws.Range("C31").Value = 9
$ws.Range("C33").Value = 1
$ws.Range("C45").Value = 24
